# Updates cryptos list price (column D) and 1h volume % (column E) values
# for the rows whose source data changed, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.207.08"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.913.56"
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3931"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09368"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.144"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.05"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.415"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "1.913.79"
$ws.Range("E14").Value = "  +2.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.333"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001129"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06626"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.236"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").Value = "28.256.69"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  +2.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.330"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.602"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.91%  "
$ws.Range("D27").Value = "2.137.23"
$ws.Range("E27").Value = "  +2.42%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.110"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1073"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.676"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.614"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.705"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06716"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02443"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2219"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.250"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("E40").Value = "  +8.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6553"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6140"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.304"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.724"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.032"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.191"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "
